# Update cryptocurrency price/volume data in the worksheet.
# Applies cell-level text updates: refreshed prices and percentage
# changes for each coin row, plus a content swap between the
# ApeXProtocol and Stellar rows (rows 48-49).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as plain text, preserving formats such as
# trailing zeros / leading zeros that Excel would otherwise mangle
# by auto-converting numeric-looking strings into real numbers.
function Set-TextValue {
    param($range, [string]$text)
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

$ws.Range('D2').Value = '66.904.46'
$ws.Range('D3').Value = '3.642.72'
$ws.Range('E3').Value = '  -8.75%  '
Set-TextValue $ws.Range('D4') '1.00'
$ws.Range('E4').Value = '  +0.02%  '
Set-TextValue $ws.Range('D5') '559.69'
$ws.Range('E5').Value = '  -8.35%  '
Set-TextValue $ws.Range('D6') '169.70'
$ws.Range('E6').Value = '  -2.11%  '
$ws.Range('D7').Value = '3.636.27'
$ws.Range('E7').Value = '  -8.79%  '
Set-TextValue $ws.Range('D8') '0.615'
$ws.Range('E8').Value = '  -10.73%  '
Set-TextValue $ws.Range('D9') '0.999'
$ws.Range('E9').Value = '  -0.05%  '
$ws.Range('E10').Value = '  -13.71%  '
$ws.Range('E11').Value = '  -16.46%  '
Set-TextValue $ws.Range('D12') '49.58'
$ws.Range('E12').Value = '  -14.03%  '
Set-TextValue $ws.Range('D13') '0.0000285'
$ws.Range('E13').Value = '  -16.13%  '
Set-TextValue $ws.Range('D14') '10.26'
$ws.Range('E14').Value = '  -12.51%  '
$ws.Range('D15').Value = '4.222.21'
$ws.Range('E15').Value = '  -8.85%  '
$ws.Range('D16').Value = '3.655.46'
$ws.Range('E16').Value = '  -8.32%  '
$ws.Range('E17').Value = '  -3.78%  '
Set-TextValue $ws.Range('D18') '19.02'
$ws.Range('E18').Value = '  -9.49%  '
$ws.Range('E19').Value = '  -12.72%  '
$ws.Range('E20').Value = '  -11.97%  '
$ws.Range('D21').Value = '66.827.06'
$ws.Range('E21').Value = '  -8.86%  '
Set-TextValue $ws.Range('D22') '397.78'
$ws.Range('E22').Value = '  -14.42%  '
Set-TextValue $ws.Range('D23') '4.35'
$ws.Range('E23').Value = '  -10.29%  '
$ws.Range('E24').Value = '  -10.94%  '
$ws.Range('E25').Value = '  -13.23%  '
Set-TextValue $ws.Range('D26') '12.47'
$ws.Range('E26').Value = '  -13.16%  '
Set-TextValue $ws.Range('D27') '10.41'
$ws.Range('E27').Value = '  -7.59%  '
Set-TextValue $ws.Range('D28') '5.97'
$ws.Range('E28').Value = '  +0.24%  '
Set-TextValue $ws.Range('D29') '3.70'
$ws.Range('E29').Value = '  -12.77%  '
Set-TextValue $ws.Range('D30') '9.20'
$ws.Range('E30').Value = '  -14.17%  '
Set-TextValue $ws.Range('D31') '31.98'
$ws.Range('E31').Value = '  -12.59%  '
Set-TextValue $ws.Range('D32') '7.47'
$ws.Range('E32').Value = '  -7.08%  '
Set-TextValue $ws.Range('D33') '12.25'
$ws.Range('E33').Value = '  -12.88%  '
Set-TextValue $ws.Range('D34') '64.00'
$ws.Range('E34').Value = '  -9.15%  '
Set-TextValue $ws.Range('D35') '0.114'
$ws.Range('E35').Value = '  -12.73%  '
Set-TextValue $ws.Range('D36') '42.04'
$ws.Range('E36').Value = '  -15.49%  '
Set-TextValue $ws.Range('D37') '580.69'
$ws.Range('E37').Value = '  -9.47%  '
$ws.Range('D38').Value = '0.0₃0876'
$ws.Range('E38').Value = '  -15.99%  '
Set-TextValue $ws.Range('D39') '0.999'
$ws.Range('E39').Value = '  -0.15%  '
Set-TextValue $ws.Range('D40') '0.999'
$ws.Range('E40').Value = '  -0.15%  '
$ws.Range('E41').Value = '  -11.13%  '
Set-TextValue $ws.Range('D42') '0.131'
$ws.Range('E42').Value = '  -12.12%  '
Set-TextValue $ws.Range('D43') '2.93'
$ws.Range('E43').Value = '  -9.64%  '
$ws.Range('E44').Value = '  -14.69%  '
Set-TextValue $ws.Range('D45') '0.0428'
$ws.Range('E45').Value = '  -12.43%  '
Set-TextValue $ws.Range('D46') '2.51'
$ws.Range('E46').Value = '  -4.67%  '
Set-TextValue $ws.Range('D47') '8.95'
$ws.Range('E47').Value = '  -16.07%  '
$ws.Range('B48').Value = 'ApeXProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
Set-TextValue $ws.Range('D48') '3.14'
$ws.Range('E48').Value = '  -8.23%  '
$ws.Range('B49').Value = 'Stellar'
$ws.Range('C49').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue $ws.Range('D49') '0.132'
$ws.Range('E49').Value = '  -12.32%  '
Set-TextValue $ws.Range('D50') '2.64'
$ws.Range('E50').Value = '  -5.41%  '
$ws.Range('D51').Value = '2.680.35'
$ws.Range('E51').Value = '  -5.03%  '
